# Updated cryptos list: refresh Price (D) and Volume/1h change (E) columns
# to match the latest scrape, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even when it
# looks numeric (e.g. "593.90"), without leaving a lasting number
# format / style change behind on the cell.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '63.363.85'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.236.93'
$ws.Range("E3").Value = '  +2.94%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue "D5" '593.90'
$ws.Range("E5").Value = '  -1.44%  '
Set-TextValue "D6" '141.48'
$ws.Range("E6").Value = '  -1.26%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.234.01'
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("E10").Value = '  -1.16%  '
Set-TextValue "D11" '5.33'
$ws.Range("E12").Value = '  -0.37%  '
Set-TextValue "D13" '0.0000246'
$ws.Range("E13").Value = '  -2.78%  '
Set-TextValue "D14" '34.30'
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("D15").Value = '3.770.44'
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = '3.239.01'
$ws.Range("E17").Value = '  +2.95%  '
$ws.Range("D18").Value = '63.373.14'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("E19").Value = '  -1.33%  '
Set-TextValue "D20" '475.32'
$ws.Range("E20").Value = '  -2.54%  '
Set-TextValue "D21" '14.10'
$ws.Range("E21").Value = '  -4.10%  '
Set-TextValue "D22" '0.723'
$ws.Range("E22").Value = '  +1.82%  '
Set-TextValue "D23" '7.89'
$ws.Range("E23").Value = '  +2.75%  '
Set-TextValue "D24" '83.85'
$ws.Range("E24").Value = '  -4.99%  '
Set-TextValue "D25" '13.16'
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("E26").Value = '  -0.01%  '
Set-TextValue "D27" '7.47'
$ws.Range("E27").Value = '  +7.09%  '
$ws.Range("E28").Value = '  -1.27%  '
Set-TextValue "D29" '8.03'
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("E30").Value = '  +2.14%  '
Set-TextValue "D31" '27.44'
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("E34").Value = '  -4.92%  '
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("E36").Value = '  -2.48%  '
Set-TextValue "D37" '52.67'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '0.0₃0712'
$ws.Range("E38").Value = '  -3.91%  '
$ws.Range("E39").Value = '  -1.07%  '
Set-TextValue "D40" '421.70'
$ws.Range("E40").Value = '  -2.82%  '
Set-TextValue "D41" '8.36'
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("D42").Value = '2.973.98'
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("E43").Value = '  -7.37%  '
$ws.Range("E44").Value = '  -7.99%  '
$ws.Range("E45").Value = '  +2.48%  '
Set-TextValue "D46" '2.15'
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  -2.83%  '
Set-TextValue "D49" '25.86'
$ws.Range("E50").Value = '  +0.02%  '
Set-TextValue "D51" '121.31'
$ws.Range("E51").Value = '  +0.72%  '
